$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price/volume table cells to the latest scraped values.
# A leading apostrophe forces a handful of decimal-looking price values (e.g. "1.010")
# to stay as literal text instead of being normalized into floating point numbers by Excel.

$ws.Range('D2').Value = "26.190.22"
$ws.Range('E2').Value = "  -4.15%  "
$ws.Range('D3').Value = "1.656.28"
$ws.Range('E3').Value = "  -3.30%  "
$ws.Range('D4').Value = "'1.010"
$ws.Range('E4').Value = "  +0.19%  "
$ws.Range('D5').Value = "'216.28"
$ws.Range('E5').Value = "  -3.66%  "
$ws.Range('D6').Value = "'0.5128"
$ws.Range('E6').Value = "  -2.72%  "
$ws.Range('D7').Value = "'1.011"
$ws.Range('E7').Value = "  +0.28%  "
$ws.Range('D8').Value = "'0.2597"
$ws.Range('D9').Value = "'0.06448"
$ws.Range('E9').Value = "  -3.05%  "
$ws.Range('D10').Value = "'19.77"
$ws.Range('E10').Value = "  -4.64%  "
$ws.Range('D11').Value = "'0.07814"
$ws.Range('E11').Value = "  +0.86%  "
$ws.Range('D12').Value = "1.660.60"
$ws.Range('E12').Value = "  -3.43%  "
$ws.Range('D13').Value = "'4.292"
$ws.Range('E13').Value = "  -3.61%  "
$ws.Range('D14').Value = "1.884.56"
$ws.Range('E14').Value = "  -3.26%  "
$ws.Range('D15').Value = "'0.5519"
$ws.Range('E15').Value = "  -4.60%  "
$ws.Range('D16').Value = "0.0₅8007"
$ws.Range('E16').Value = "  -1.86%  "
$ws.Range('E17').Value = "  -5.19%  "
$ws.Range('D18').Value = "26.219.39"
$ws.Range('E18').Value = "  -4.05%  "
$ws.Range('E19').Value = "  +0.18%  "
$ws.Range('D20').Value = "'209.17"
$ws.Range('E20').Value = "  -4.51%  "
$ws.Range('D21').Value = "'4.404"
$ws.Range('E21').Value = "  -5.00%  "
$ws.Range('D22').Value = "'10.09"
$ws.Range('E22').Value = "  -2.90%  "
$ws.Range('D23').Value = "'6.064"
$ws.Range('E23').Value = "  +0.80%  "
$ws.Range('E24').Value = "  +0.32%  "
$ws.Range('D25').Value = "'1.834"
$ws.Range('E25').Value = "  +7.44%  "
$ws.Range('D26').Value = "'144.42"
$ws.Range('E26').Value = "  -0.58%  "
$ws.Range('E27').Value = "  -2.46%  "
$ws.Range('D28').Value = "'6.970"
$ws.Range('D29').Value = "'15.84"
$ws.Range('E29').Value = "  -1.84%  "
$ws.Range('D30').Value = "'0.05078"
$ws.Range('E30').Value = "  -5.32%  "
$ws.Range('E31').Value = "  -3.80%  "
$ws.Range('D32').Value = "'3.351"
$ws.Range('E32').Value = "  -3.53%  "
$ws.Range('E33').Value = "  -3.99%  "
$ws.Range('D34').Value = "'1.555"
$ws.Range('E34').Value = "  -4.67%  "
$ws.Range('E35').Value = "  -3.77%  "
$ws.Range('D36').Value = "'2.359"
$ws.Range('E36').Value = "  -1.63%  "
$ws.Range('D37').Value = "'0.9219"
$ws.Range('E37').Value = "  -2.97%  "
$ws.Range('D38').Value = "1.173.28"
$ws.Range('E38').Value = "  +1.65%  "
$ws.Range('D39').Value = "'0.5709"
$ws.Range('E39').Value = "  -2.62%  "
$ws.Range('D40').Value = "'0.01586"
$ws.Range('E40').Value = "  -3.61%  "
$ws.Range('B41').Value = "PaxDollar"
$ws.Range('C41').Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range('D41').Value = "'1.011"
$ws.Range('E41').Value = "  +0.28%  "
$ws.Range('B42').Value = "mCoin"
$ws.Range('C42').Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range('D42').Value = "'2.566"
$ws.Range('E42').Value = "  -0.35%  "
$ws.Range('E43').Value = "  -2.57%  "
$ws.Range('D44').Value = "'0.8268"
$ws.Range('E44').Value = "  -1.52%  "
$ws.Range('D45').Value = "'100.57"
$ws.Range('E45').Value = "  -0.44%  "
$ws.Range('D46').Value = "1.795.79"
$ws.Range('E46').Value = "  -3.21%  "
$ws.Range('D47').Value = "0.0₈112"
$ws.Range('E47').Value = "  -4.50%  "
$ws.Range('E48').Value = "  +0.13%  "
$ws.Range('D49').Value = "'55.44"
$ws.Range('E49').Value = "  -3.37%  "
$ws.Range('D50').Value = "'1.004"
$ws.Range('E50').Value = "  -0.20%  "
$ws.Range('D51').Value = "'7.872"
$ws.Range('E51').Value = "  -3.24%  "
